# Congo King Quad Shot doc edit:
#  1. Insert a new "Meta description" paragraph right after the title (Heading1).
#  2. Remove the duplicated bold title paragraph near the end of the document.
#  3. Replace the italic meta-description paragraph's text with the DALLE image prompt.

$d = $word.ActiveDocument

# --- 1. Insert "Meta description" paragraph after the title paragraph ---
$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range
$titleRange.Collapse(0)  # wdCollapseEnd
$titleRange.InsertParagraphAfter() | Out-Null

$metaRange = $d.Paragraphs(2).Range
$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Experience the thrill of winning with Congo King Quad Shot. 100 paylines, x4 jackpots, and 15 free spins. Play free online here.</w:t></w:r></w:p>'
$metaRange.InsertXML($metaXml) | Out-Null

# --- 2. Delete the bold "Play Congo King Quad Shot..." paragraph that used to sit ---
#        right before the final (italic) paragraph. Search from the end backwards so
#        we find the duplicated copy rather than the Heading1 title at the top.
$oldTitle = "Play Congo King Quad Shot Free Online Review | 100 Paylines"
$lastIndex = $d.Paragraphs.Count
for ($i = $lastIndex; $i -ge 1; $i--) {
    $para = $d.Paragraphs($i)
    $text = $para.Range.Text.Trim()
    if ($text -eq $oldTitle) {
        $para.Range.Delete() | Out-Null
        break
    }
}

# --- 3. Replace the remaining italic description text with the DALLE prompt ---
$oldDescription = "Experience the thrill of winning with Congo King Quad Shot. 100 paylines, x4 jackpots, and 15 free spins. Play free online here."
$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>DALLE, please create a feature image fitting the game "Congo King Quad Shot" that meets the following requirements: - The image should be in cartoon style - The image should feature a happy Maya warrior with glasses. The image should capture the adventurous spirit of the game and convey the excitement of exploring through the jungle. Please ensure the colors used in the image are vibrant and eye-catching. The image should be appealing and encourage players to take a chance on the game.</w:t></w:r></w:p>'

$n = $d.Paragraphs.Count
for ($i = $n; $i -ge 1; $i--) {
    $para = $d.Paragraphs($i)
    $text = $para.Range.Text.Trim()
    if ($text -eq $oldDescription) {
        $para.Range.InsertXML($newXml) | Out-Null
        break
    }
}

Write-Output "Edit complete."
